$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: caster wheel now sourced from Home Depot, price confirmed, spec note added ---
$ws.Range("B14").Value = "Home Depot"
$ws.Range("C14").Value = 2.97
$ws.Range("E14").Value = "40 lb load rating"

# --- Row 15: new BOM line for the drive wheel ---
$ws.Range("A15").Value = "drive wheel"
$ws.Range("B15").Value = "Home Depot"
$ws.Range("C15").Value = 8.14
$ws.Range("C15").Style = "Currency"
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Formula = "=D15*F15"

# Link the new part name to its source, like the other BOM rows
$ws.Hyperlinks.Add($ws.Range("A15"), "https://www.homedepot.com/p/Shepherd-3-in-Non-Marking-Rubber-Flat-Free-Wheel-with-1-2-in-Diameter-Ball-Bearing-Axle-9731/202526123")

# Adding the hyperlink re-stamps A15 with the generic "Hyperlink" style; restore the
# same look the rest of column A already uses (it's styled as a hyperlink by default)
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)

# Reflect the final reviewed selection
$ws.Range("I2").Select()

$wb.Save()
